# "Generate Report for Handback"
#
# For the e0857031-8132-4245-b047-bff43da63f9b file row (row 7) on both the
# zh-cn and de-de report sheets, the handback tool discovered a target file
# whose commit (a80d33aa...) is not the current "latest" commit
# (ace668b7...). It now fills in the "Latest Target File" / "Latest Handback
# File" / "Latest Handback DateTime" columns and records a warning in
# "Error Detail". The "Error Detail" column is also widened so the message
# is readable.

$wb = $excel.ActiveWorkbook

$fileMd       = "e0857031-8132-4245-b047-bff43da63f9b.md"
$currentRev   = "a80d33aafeb6b8a0405122afb8551bb4a799ca56"
$latestRev    = "ace668b759f41b6eb51a84f1d625d1fae1cc9c58"
$targetUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$currentRev/e2e/$fileMd"
$latestUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$latestRev/e2e/$fileMd"
$errorDetail  = "The version of handback file is not the latest, current: $targetUrl, latest: $latestUrl."

# cornflowerblue (FF6495ED), matching the workbook's existing "HyperLink" cell style
$hyperlinkColor = 15570276

$sheets = @{
    "zh-cn" = @{ XlfName = "e0857031-8132-4245-b047-bff43da63f9b.a58e485f01871edbf7439d5deac0b7242850f4d9.zh-cn.xlf"; HandbackDate = "2016-08-24 16:44:47" }
    "de-de" = @{ XlfName = "e0857031-8132-4245-b047-bff43da63f9b.a58e485f01871edbf7439d5deac0b7242850f4d9.de-de.xlf"; HandbackDate = "2016-08-24 16:44:54" }
}

foreach ($sheetName in $sheets.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $info = $sheets[$sheetName]

    # I7 - Latest Target File: hyperlink to the (non-latest) target md file
    $ws.Hyperlinks.Add($ws.Range("I7"), $targetUrl, $null, "", $fileMd)
    $ws.Range("I7").Font.Color = $hyperlinkColor
    $ws.Range("I7").Font.Underline = $true

    # J7 - Latest Handback File
    $ws.Range("J7").Value = $info.XlfName

    # K7 - Latest Handback DateTime
    $ws.Range("K7").Value = $info.HandbackDate

    # P7 - Error Detail
    $ws.Range("P7").Value = $errorDetail

    # Widen the Error Detail column (P / 16) so the message is readable
    $ws.Range("P1").ColumnWidth = 39.1667
}
